$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "JD_003"
$ws.Range("B4").Value = "Cyber Security Engineer"
$ws.Range("C4").Value = "We are seeking a Cyber Security Engineer  to build and maintain high-quality software solutions.`nWork with global teams to drive innovation and deliver scalable applications.`nJoin Akkodis and be part of a tech-driven, collaborative environment."
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 5

$ws.Rows.Item(4).EntireRow.AutoFit()
